$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 101
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(19, 6).Value = 3869
$ws.Cells.Item(20, 6).Value = 6200
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 2546
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 296
$ws.Cells.Item(35, 6).Value = 365
$ws.Cells.Item(36, 6).Value = 163
$ws.Cells.Item(38, 6).Value = 944
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 107

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 4766
$ws.Cells.Item(5, 6).Value = 207
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 737
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 141
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(30, 6).Value = 2546
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 365
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(46, 6).Value = 0
